# Generate Report for Handoff
# - Update the status / handoff timestamp for the remaining (14471e81...) row
#   from "Handed back: in sync with en-US" -> "Ready for handoff" with new
#   handoff datetimes.
# - Remove the second file's row (ef3b9127-46ba-4622-b04d-80d8965d3e01...)
#   entirely from every sheet (Overview, zh-cn, de-de), including its
#   hyperlinks.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet 1: Overview  (columns: A=File Name, B=zh-cn, C=de-de, D=Latest Handoff Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(2, 2).Value2 = $newStatus
$wsOverview.Cells.Item(2, 3).Value2 = $newStatus
$wsOverview.Cells.Item(2, 4).Value2 = "2016-03-22 10:51:24"

# Remove row 3 (ef3b9127-...) along with its hyperlink, then delete the row.
$wsOverview.Range("A3").Hyperlinks.Delete()
$wsOverview.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# (columns: A=Source File Name, B=File Extension, C=Status,
#  D=Latest Handoff File, E=Latest Handoff Datetime, F=Latest Target File,
#  G=Latest Handback File, H=Latest Handback DateTime, I=Reference Tokens,
#  J=Handoff Reason, K=Dependency From, L=Error Detail)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(2, 3).Value2 = $newStatus
$wsZhCn.Cells.Item(2, 5).Value2 = "2016-03-22 10:51:17"

# Remove row 3 (ef3b9127-...) along with its hyperlinks, then delete the row.
$wsZhCn.Range("A3").Hyperlinks.Delete()
$wsZhCn.Range("D3").Hyperlinks.Delete()
$wsZhCn.Range("F3").Hyperlinks.Delete()
$wsZhCn.Range("G3").Hyperlinks.Delete()
$wsZhCn.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# Sheet 3: de-de (same column layout as zh-cn)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(2, 3).Value2 = $newStatus
$wsDeDe.Cells.Item(2, 5).Value2 = "2016-03-22 10:51:24"

# Remove row 3 (ef3b9127-...) along with its hyperlinks, then delete the row.
$wsDeDe.Range("A3").Hyperlinks.Delete()
$wsDeDe.Range("D3").Hyperlinks.Delete()
$wsDeDe.Range("F3").Hyperlinks.Delete()
$wsDeDe.Range("G3").Hyperlinks.Delete()
$wsDeDe.Rows.Item(3).Delete()
